# feat: add 2022-Q1 data
#
# The workbook currently has 3 sheets: 2021-Q3, 2021-Q4, 总计 (summary).
# We turn the old 总计 sheet into the new 2022-Q1 holdings sheet (it keeps
# sheetId 3 / rId3, matching how Excel would reuse the tab when it's
# repurposed), and append a brand new 总计 sheet after it that aggregates
# all four quarters now available.

$wb = $excel.ActiveWorkbook
$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Grab a still-pristine sheet to borrow cell styles/formatting from (the
# bold/boxed header style and the index-column style used throughout the
# workbook).
$ref = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) Repurpose the existing "总计" sheet (3rd tab) into "2022-Q1"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# Header row (B1:H1) -- reuse the bold/boxed header style used elsewhere
# in the workbook across the whole header span.
$ref.Range("B1").Copy() | Out-Null
$q1.Range("B1:H1").PasteSpecial($fmt) | Out-Null

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Index column (A2:A5) -- reuse the centered/boxed index style.
$ref.Range("A2").Copy() | Out-Null
$q1.Range("A2:A5").PasteSpecial($fmt) | Out-Null

# Index column + text-valued data columns need to hold numeric-looking
# strings ("012751", "0.34", ...) verbatim -- format as Text first so
# Excel doesn't silently coerce them to numbers and drop leading zeros.
$q1.Range("B2:G5").NumberFormat = "@"

$q1Data = @(
    @(0, "012751", "建信纳斯达克100指数（QDII）A 美元现汇", "0.34", "88.02", "6.09", "0.0207", 4),
    @(1, "012752", "建信纳斯达克100指数（QDII）C 人民币",   "0.34", "88.02", "6.09", "0.0207", 4),
    @(2, "012753", "建信纳斯达克100指数（QDII）C 美元现汇", "0.34", "88.02", "6.09", "0.0207", 4),
    @(3, "539002", "建信新兴市场优选混合QDII",               "0.14", "83.76", "3.82", "0.0053", 9)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Range("A$r").Value = $row[0]
    $q1.Range("B$r").Value = $row[1]
    $q1.Range("C$r").Value = $row[2]
    $q1.Range("D$r").Value = $row[3]
    $q1.Range("E$r").Value = $row[4]
    $q1.Range("F$r").Value = $row[5]
    $q1.Range("G$r").Value = $row[6]
    $q1.Range("H$r").Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------
# 2) Insert a fresh "总计" sheet right after "2022-Q1"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

# Match the page-margin settings the rest of the workbook's "summary"
# sheets use (0.75in / 1in / 0.5in, expressed in points for the COM API).
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$ref.Range("B1").Copy() | Out-Null
$total.Range("B1:D1").PasteSpecial($fmt) | Out-Null

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$ref.Range("A2").Copy() | Out-Null
$total.Range("A2:A4").PasteSpecial($fmt) | Out-Null

$totalData = @(
    @(0, "2022-Q1", 4, 0.07000000000000001),
    @(1, "2021-Q4", 3, 1.2),
    @(2, "2021-Q3", 4, 0.02)
)

$r = 2
foreach ($row in $totalData) {
    $total.Range("A$r").Value = $row[0]
    $total.Range("B$r").Value = $row[1]
    $total.Range("C$r").Value = $row[2]
    $total.Range("D$r").Value = $row[3]
    $r++
}

Write-Host "done"
